# EIA Table A.3.B update: 2017-01-31 update (EPM 2016_11 data), chunk 7
# - Update subtitle from "October 2016" to "November 2016"
# - Update 129 data cell values (Relative Standard Error percentages) for rows 4-63

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subtitle text in A2
$ws.Range("A2").Value = "Independent Power Producers by Census Division and State, Year-to-Date through November 2016"

# Update data cell values (columns B-H, rows 4-63)
$ws.Range("C4").Value = 4
$ws.Range("E4").Value = 2
$ws.Range("H4").Value = 16
$ws.Range("C5").Value = 69
$ws.Range("H5").Value = 97
$ws.Range("C6").Value = 173
$ws.Range("E6").Value = 2
$ws.Range("H6").Value = 20
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = 5
$ws.Range("H7").Value = 38
$ws.Range("C8").Value = 22
$ws.Range("H8").Value = 33
$ws.Range("C9").Value = 29
$ws.Range("E9").Value = 0
$ws.Range("H9").Value = 989
$ws.Range("H10").Value = 55
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 42
$ws.Range("F11").Value = 79
$ws.Range("H11").Value = 16
$ws.Range("C12").Value = 229
$ws.Range("E12").Value = 2
$ws.Range("H12").Value = 886
$ws.Range("C13").Value = 142
$ws.Range("E13").Value = 3
$ws.Range("H13").Value = 24
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 28
$ws.Range("F14").Value = 79
$ws.Range("H14").Value = 17
$ws.Range("F15").Value = 13
$ws.Range("H15").Value = 53
$ws.Range("H16").Value = 65
$ws.Range("E17").Value = 5
$ws.Range("H18").Value = 119
$ws.Range("C19").Value = 2
$ws.Range("F19").Value = 59
$ws.Range("H19").Value = 97
$ws.Range("E20").Value = 0.18
$ws.Range("H20").Value = 111
$ws.Range("B21").Value = 186
$ws.Range("C21").Value = 438
$ws.Range("E21").Value = 7
$ws.Range("H21").Value = 76
$ws.Range("C22").Value = 120
$ws.Range("E22").Value = 12508
$ws.Range("H22").Value = 449
$ws.Range("H23").Value = 354
$ws.Range("C24").Value = 500
$ws.Range("E24").Value = 11
$ws.Range("H24").Value = 79
$ws.Range("B25").Value = 186
$ws.Range("C25").Value = 2159
$ws.Range("E25").Value = 9
$ws.Range("C26").Value = 594
$ws.Range("C27").Value = 34
$ws.Range("E27").Value = 3
$ws.Range("H27").Value = 10
$ws.Range("C28").Value = 431
$ws.Range("E28").Value = 8
$ws.Range("C30").Value = 343
$ws.Range("E30").Value = 9
$ws.Range("C31").Value = 113
$ws.Range("E31").Value = 4
$ws.Range("H31").Value = 506
$ws.Range("C32").Value = 32
$ws.Range("E32").Value = 20
$ws.Range("H32").Value = 0
$ws.Range("B33").Value = 250
$ws.Range("C33").Value = 28
$ws.Range("H33").Value = 248
$ws.Range("C34").Value = 464
$ws.Range("E34").Value = 3
$ws.Range("H34").Value = 201
$ws.Range("C35").Value = 87
$ws.Range("H35").Value = 172
$ws.Range("E36").Value = 9
$ws.Range("H36").Value = 12
$ws.Range("C37").Value = 122
$ws.Range("E37").Value = 0
$ws.Range("H37").Value = 500
$ws.Range("C38").Value = 124
$ws.Range("E38").Value = 0
$ws.Range("H39").Value = 500
$ws.Range("C41").Value = 774
$ws.Range("E42").Value = 1
$ws.Range("H42").Value = 13
$ws.Range("H43").Value = 217
$ws.Range("E44").Value = 0.38
$ws.Range("E46").Value = 1
$ws.Range("H46").Value = 237
$ws.Range("C47").Value = 23
$ws.Range("E47").Value = 2
$ws.Range("H47").Value = 39
$ws.Range("E48").Value = 0
$ws.Range("B49").Value = 138
$ws.Range("E49").Value = 3
$ws.Range("H49").Value = 83
$ws.Range("E50").Value = 9
$ws.Range("H50").Value = 52
$ws.Range("E51").Value = 180
$ws.Range("H51").Value = 95
$ws.Range("E52").Value = 6
$ws.Range("H52").Value = 289
$ws.Range("E53").Value = 3
$ws.Range("B54").Value = 90
$ws.Range("C54").Value = 485
$ws.Range("E54").Value = 102
$ws.Range("H54").Value = 382
$ws.Range("B55").Value = 89
$ws.Range("E55").Value = 1188
$ws.Range("H55").Value = 359
$ws.Range("C56").Value = 63
$ws.Range("E56").Value = 2
$ws.Range("H56").Value = 30
$ws.Range("C57").Value = 232
$ws.Range("H57").Value = 43
$ws.Range("H58").Value = 64
$ws.Range("C59").Value = 42
$ws.Range("H59").Value = 50
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = 16
$ws.Range("B61").Value = 65
$ws.Range("C62").Value = 16
$ws.Range("C63").Value = 7
$ws.Range("E63").Value = 0.47
$ws.Range("F63").Value = 7
$ws.Range("H63").Value = 9
